$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the column header in A1 from "Goal Name" to "Goal"
$ws.Range("A1").Value = "Goal"

# Update the active selection to A3 (as recorded in the saved view state)
$ws.Range("A3").Select()
